# Kosaraju_LabExam03Grading.xlsx - grading update
# Awards full "Total Points" (10/10) for two CustomerMapping-class rubric
# items that were previously left blank, and updates the sheet's active
# selection/scroll position to reflect where the grader was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
[void]$ws.Activate()

# Question 14 "whoPurchasedProduct() method" - full credit (10/10)
$ws.Range("E22").Value = 10

# Question 16 "findAllBrands()" - full credit (10/10)
$ws.Range("E24").Value = 10

# (Totals in E26/E38 are formulas - SUM(E18:E25) / SUM(E7,E15,E26,E35,E37,E31) -
#  Excel recalculates them automatically from the two edits above.)

# Reflect the grader's final scroll position / selection on the sheet.
try {
    $excel.ActiveWindow.ScrollRow = 11
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Older/limited hosts may not expose ActiveWindow scrolling; selection
    # below still captures the important state.
}
$ws.Range("F23").Select() | Out-Null
